# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Values are written as text (matching the sheet's existing inline-string cells);
# price cells whose new text parses as a plain number get NumberFormat "@" set
# first so Excel keeps them as text instead of silently recasting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.254.69'
$ws.Range("E2").Value = '  +0.55%  '

# Row 3
$ws.Range("D3").Value = '1.590.83'
$ws.Range("E3").Value = '  +1.23%  '

# Row 4
$ws.Range("E4").Value = '  -0.26%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.36'
$ws.Range("E5").Value = '  +1.86%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.502'
$ws.Range("E6").Value = '  +1.19%  '

# Row 7
$ws.Range("E7").Value = '  -0.22%  '

# Row 8
$ws.Range("E8").Value = '  +0.68%  '

# Row 9
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.41'
$ws.Range("E10").Value = '  -0.49%  '

# Row 11
$ws.Range("E11").Value = '  +0.59%  '

# Row 12
$ws.Range("D12").Value = '1.812.97'
$ws.Range("E12").Value = '  +1.09%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.597.20'
$ws.Range("E13").Value = '  +1.92%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("E14").Value = '  +0.03%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("E15").Value = '  +1.96%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.43'
$ws.Range("E16").Value = '  +0.30%  '

# Row 17
$ws.Range("D17").Value = '26.257.14'
$ws.Range("E17").Value = '  +0.53%  '

# Row 18
$ws.Range("E18").Value = '  +0.30%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.44'
$ws.Range("E19").Value = '  +2.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '213.97'
$ws.Range("E20").Value = '  +3.58%  '

# Row 21
$ws.Range("E21").Value = '  -0.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.29'
$ws.Range("E22").Value = '  +1.14%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.01'
$ws.Range("E23").Value = '  +1.94%  '

# Row 24
$ws.Range("E24").Value = '  -2.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.66'
$ws.Range("E25").Value = '  -0.33%  '

# Row 26
$ws.Range("E26").Value = '  -0.24%  '

# Row 27
$ws.Range("E27").Value = '  +1.60%  '

# Row 28
$ws.Range("E28").Value = '  +0.20%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.21'
$ws.Range("E29").Value = '  +0.12%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0499'
$ws.Range("E30").Value = '  -1.25%  '

# Row 31
$ws.Range("E31").Value = '  +1.36%  '

# Row 32
$ws.Range("E32").Value = '  +0.16%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.95'
$ws.Range("E33").Value = '  -0.73%  '

# Row 34
$ws.Range("D34").Value = '1.336.53'
$ws.Range("E34").Value = '  +4.56%  '

# Row 35
$ws.Range("E35").Value = '  -0.62%  '

# Row 36
$ws.Range("E36").Value = '  +0.02%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.591'
$ws.Range("E37").Value = '  -3.14%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0167'
$ws.Range("E38").Value = '  +0.81%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.817'
$ws.Range("E39").Value = '  +0.67%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.75'
$ws.Range("E40").Value = '  +3.36%  '

# Row 41
$ws.Range("E41").Value = '  -0.20%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("E42").Value = '  -7.13%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.15'
$ws.Range("E43").Value = '  +0.29%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.768'
$ws.Range("E44").Value = '  +0.85%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.87'
$ws.Range("E45").Value = '  -0.58%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.724.89'
$ws.Range("E46").Value = '  +0.98%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.64'
$ws.Range("E47").Value = '  -3.84%  '

# Row 48
$ws.Range("E48").Value = '  -1.85%  '

# Row 49
$ws.Range("E49").Value = '  -0.73%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0974'
$ws.Range("E50").Value = '  -2.67%  '

# Row 51
$ws.Range("E51").Value = '  -0.39%  '
